$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("prep_sheet") ---------------------------------------------

# Switch margins to the "Narrow" preset (in points; COM margins are always points)
$ws1.PageSetup.LeftMargin = $excel.InchesToPoints(0.25)
$ws1.PageSetup.RightMargin = $excel.InchesToPoints(0.25)
$ws1.PageSetup.TopMargin = $excel.InchesToPoints(0.75)
$ws1.PageSetup.BottomMargin = $excel.InchesToPoints(0.75)
$ws1.PageSetup.HeaderMargin = $excel.InchesToPoints(0.3)
$ws1.PageSetup.FooterMargin = $excel.InchesToPoints(0.3)

# Print scale 73% and keep portrait orientation
$ws1.PageSetup.Zoom = 73
$ws1.PageSetup.Orientation = 1

# sheetPr/pageSetUpPr fitToPage="1"  (enable "fit to page" print scaling mode)
$ws1.PageSetup.FitToPagesWide = 1
$ws1.PageSetup.FitToPagesTall = 1

# Move the view / selection: scroll so row 4 is the top visible row, and
# select H53
$ws1.Activate()
$ws1.Range("H53").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# --- Sheet2 ("order_sheet") ---------------------------------------------

$ws2.PageSetup.Orientation = 1
